{"js": "// Insert a new paragraph with the text \"Change !!!!!!!!\" right after the\n// first paragraph (\"New file \") and before the final (empty) paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertParagraph(\"Change !!!!!!!!\", \"After\");\nawait context.sync();\n", "ps1": "# Insert a new paragraph with the text \"Change !!!!!!!!\" right after the\n# first paragraph (\"New file \") and before the final (empty) paragraph.\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs.Item(1)\n$firstParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item(2)\n$newParagraph.Range.InsertAfter(\"Change !!!!!!!!\")\n"}
